$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A14").Value = 45978
$ws.Range("B14").Value = 3

$ws.Range("A14:B14").Select()
